$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033286581860451
$ws.Range("D2").Value = 1.038334395248997
$ws.Range("E2").Value = 1.043063984897162
$ws.Range("F2").Value = 1.05512392456556
$ws.Range("I2").Value = 1.040041983958494
$ws.Range("J2").Value = 1.03841213996228
$ws.Range("K2").Value = 1.041122723855539
$ws.Range("L2").Value = 1.045838905731709
$ws.Range("M2").Value = 1.057865246688045
$ws.Range("N2").Value = 1.016782337452966
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034085232048705
$ws.Range("D3").Value = 1.03894417452142
$ws.Range("E3").Value = 1.043810493469389
$ws.Range("F3").Value = 1.056036609161052
$ws.Range("I3").Value = 1.040250482773939
$ws.Range("J3").Value = 1.038854205460048
$ws.Range("K3").Value = 1.04154302438993
$ws.Range("L3").Value = 1.046396541905605
$ws.Range("M3").Value = 1.058591056900044
$ws.Range("N3").Value = 1.016929443183196
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034602544887663
$ws.Range("D4").Value = 1.039339094099636
$ws.Range("E4").Value = 1.044294433257219
$ws.Range("F4").Value = 1.05662834880711
$ws.Range("I4").Value = 1.040384326120727
$ws.Range("J4").Value = 1.039140106620779
$ws.Range("K4").Value = 1.041814643795819
$ws.Range("L4").Value = 1.046757601577994
$ws.Range("M4").Value = 1.059061241960764
$ws.Range("N4").Value = 1.017024559185061
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034820149051647
$ws.Range("D5").Value = 1.039505200678975
$ws.Range("E5").Value = 1.044498094888756
$ws.Range("F5").Value = 1.056877394723071
$ws.Range("I5").Value = 1.040440337033781
$ws.Range("J5").Value = 1.039260263573682
$ws.Range("K5").Value = 1.041928749155034
$ws.Range("L5").Value = 1.046909445116254
$ws.Range("M5").Value = 1.059259034707658
$ws.Range("N5").Value = 1.017064528451253
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034856693078654
$ws.Range("D6").Value = 1.039533095483804
$ws.Range("E6").Value = 1.044532303030093
$ws.Range("F6").Value = 1.056919226907401
$ws.Range("I6").Value = 1.04044972642779
$ws.Range("J6").Value = 1.039280436313965
$ws.Range("K6").Value = 1.041947903016524
$ws.Range("L6").Value = 1.046934943442335
$ws.Range("M6").Value = 1.05929225237538
$ws.Range("N6").Value = 1.017071238427695
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034605452032787
$ws.Range("D7").Value = 1.0393413133024
$ws.Range("E7").Value = 1.044297153757915
$ws.Range("F7").Value = 1.05663167547867
$ws.Range("I7").Value = 1.040385075551492
$ws.Range("J7").Value = 1.039141712306457
$ws.Range("K7").Value = 1.04181616880669
$ws.Range("L7").Value = 1.0467596303082
$ws.Range("M7").Value = 1.05906388438043
$ws.Range("N7").Value = 1.017025093325679
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033556378097326
$ws.Range("D8").Value = 1.038540399130587
$ws.Range("E8").Value = 1.043316083753575
$ws.Range("F8").Value = 1.055432126829514
$ws.Range("I8").Value = 1.040112668021341
$ws.Range("J8").Value = 1.038561567301311
$ws.Range("K8").Value = 1.041264836664253
$ws.Range("L8").Value = 1.046027312632356
$ws.Range("M8").Value = 1.058110425376584
$ws.Range("N8").Value = 1.016832067006731
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03171193124655
$ws.Range("D9").Value = 1.037131857972484
$ws.Range("E9").Value = 1.041594278123492
$ws.Range("F9").Value = 1.053327428142303
$ws.Range("I9").Value = 1.03962449811534
$ws.Range("J9").Value = 1.037538223470402
$ws.Range("K9").Value = 1.04029074789939
$ws.Range("L9").Value = 1.044738713995173
$ws.Range("M9").Value = 1.056434498700415
$ws.Range("N9").Value = 1.016491404020991
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030485195063617
$ws.Range("D10").Value = 1.036194800856723
$ws.Range("E10").Value = 1.040451200128216
$ws.Range("F10").Value = 1.051930502045378
$ws.Range("I10").Value = 1.039293619386841
$ws.Range("J10").Value = 1.036855359065626
$ws.Range("K10").Value = 1.039639703574967
$ws.Range("L10").Value = 1.043880970111479
$ws.Range("M10").Value = 1.05532013198781
$ws.Range("N10").Value = 1.016263968283122
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029954712708437
$ws.Range("D11").Value = 1.03578953464961
$ws.Range("E11").Value = 1.039957394575644
$ws.Range("F11").Value = 1.051327114034239
$ws.Range("I11").Value = 1.039149067234283
$ws.Range("J11").Value = 1.036559536099111
$ws.Range("K11").Value = 1.039357418610434
$ws.Range("L11").Value = 1.043509889143842
$ws.Range("M11").Value = 1.054838312181585
$ws.Range("N11").Value = 1.016165414116109
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029757775195478
$ws.Range("D12").Value = 1.035639075532657
$ws.Range("E12").Value = 1.03977414857358
$ws.Range("F12").Value = 1.051103214717097
$ws.Range("I12").Value = 1.039095182723312
$ws.Range("J12").Value = 1.036449634764609
$ws.Range("K12").Value = 1.039252509844092
$ws.Range("L12").Value = 1.043372103854385
$ws.Range("M12").Value = 1.054659450984028
$ws.Range("N12").Value = 1.016128796194071
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029800014079611
$ws.Range("D13").Value = 1.03567134611001
$ws.Range("E13").Value = 1.039813447497228
$ws.Range("F13").Value = 1.051151231628577
$ws.Range("I13").Value = 1.039106749782113
$ws.Range("J13").Value = 1.036473209836946
$ws.Range("K13").Value = 1.039275015624018
$ws.Range("L13").Value = 1.043401656948832
$ws.Range("M13").Value = 1.054697812397384
$ws.Range("N13").Value = 1.016136651331316
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029938431594347
$ws.Range("D14").Value = 1.035777096115692
$ws.Range("E14").Value = 1.03994224381567
$ws.Range("F14").Value = 1.051308601828051
$ws.Range("I14").Value = 1.039144617023647
$ws.Range("J14").Value = 1.036550452011378
$ws.Range("K14").Value = 1.03934874794501
$ws.Range("L14").Value = 1.043498498726057
$ws.Range("M14").Value = 1.054823525243443
$ws.Range("N14").Value = 1.016162387479179
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030023729456834
$ws.Range("D15").Value = 1.035842262150216
$ws.Range("E15").Value = 1.040021622772199
$ws.Range("F15").Value = 1.051405592788859
$ws.Range("I15").Value = 1.039167922942483
$ws.Range("J15").Value = 1.036598040915923
$ws.Range("K15").Value = 1.039394169537843
$ws.Range("L15").Value = 1.043558172903571
$ws.Range("M15").Value = 1.054900995456652
$ws.Range("N15").Value = 1.016178242989217
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030520416377342
$ws.Range("D16").Value = 1.036221707429314
$ws.Range("E16").Value = 1.04048399686619
$ws.Range("F16").Value = 1.051970578485489
$ws.Range("I16").Value = 1.039303185950107
$ws.Range("J16").Value = 1.036874989056838
$ws.Range("K16").Value = 1.039658430035308
$ws.Range("L16").Value = 1.043905604563793
$ws.Range("M16").Value = 1.055352123851444
$ws.Range("N16").Value = 1.016270507499169
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03083216424774
$ws.Range("D17").Value = 1.036459854754994
$ws.Range("E17").Value = 1.0407743424077
$ws.Range("F17").Value = 1.052325379446866
$ws.Range("I17").Value = 1.039387690888035
$ws.Range("J17").Value = 1.037048675195416
$ws.Range("K17").Value = 1.039824093252881
$ws.Range("L17").Value = 1.044123628053664
$ws.Range("M17").Value = 1.055635295624536
$ws.Range("N17").Value = 1.016328363364455
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031014069173191
$ws.Range("D18").Value = 1.036598808802482
$ws.Range("E18").Value = 1.040943807341863
$ws.Range("F18").Value = 1.052532472496136
$ws.Range("I18").Value = 1.039436857682188
$ws.Range("J18").Value = 1.037149969966298
$ws.Range("K18").Value = 1.039920685200329
$ws.Range("L18").Value = 1.044250828996348
$ws.Range("M18").Value = 1.055800533157983
$ws.Range("N18").Value = 1.01636210263905
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031076105477438
$ws.Range("D19").Value = 1.036646196414395
$ws.Range("E19").Value = 1.04100160932499
$ws.Range("F19").Value = 1.052603110180024
$ws.Range("I19").Value = 1.039453601314587
$ws.Range("J19").Value = 1.037184506552122
$ws.Range("K19").Value = 1.039953614309922
$ws.Range("L19").Value = 1.044294206534292
$ws.Range("M19").Value = 1.055856886407523
$ws.Range("N19").Value = 1.016373605646425
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030798709640145
$ws.Range("D20").Value = 1.036434298952887
$ws.Range("E20").Value = 1.040743179553898
$ws.Range("F20").Value = 1.052287297788781
$ws.Range("I20").Value = 1.039378637076861
$ws.Range("J20").Value = 1.037030041695105
$ws.Range("K20").Value = 1.039806322934257
$ws.Range("L20").Value = 1.044100232924776
$ws.Range("M20").Value = 1.055604906903584
$ws.Range("N20").Value = 1.016322156705977
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029897668111762
$ws.Range("D21").Value = 1.035745953295372
$ws.Range("E21").Value = 1.039904311654005
$ws.Range("F21").Value = 1.051262253976104
$ws.Range("I21").Value = 1.039133471341948
$ws.Range("J21").Value = 1.036527706644877
$ws.Range("K21").Value = 1.039327037137549
$ws.Range("L21").Value = 1.043469979824997
$ws.Range("M21").Value = 1.054786502944876
$ws.Range("N21").Value = 1.01615480911178
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029331768924013
$ws.Range("D22").Value = 1.035313596712513
$ws.Range("E22").Value = 1.039377897657263
$ws.Range("F22").Value = 1.050619075956232
$ws.Range("I22").Value = 1.038978218909505
$ws.Range("J22").Value = 1.036211756380105
$ws.Range("K22").Value = 1.039025370313996
$ws.Range("L22").Value = 1.043074008766879
$ws.Range("M22").Value = 1.054272566357356
$ws.Range("N22").Value = 1.016049530415692
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029631703319205
$ws.Range("D23").Value = 1.03554275540385
$ws.Range("E23").Value = 1.039656862788632
$ws.Range("F23").Value = 1.050959912180334
$ws.Range("I23").Value = 1.039060625817795
$ws.Range("J23").Value = 1.036379257820508
$ws.Range("K23").Value = 1.039185319581077
$ws.Range("L23").Value = 1.043283892154991
$ws.Range("M23").Value = 1.054544953912919
$ws.Range("N23").Value = 1.016105346240549
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03081382612841
$ws.Range("D24").Value = 1.036445846375636
$ws.Range("E24").Value = 1.04075726036251
$ws.Range("F24").Value = 1.052304504807252
$ws.Range("I24").Value = 1.039382728486127
$ws.Range("J24").Value = 1.037038461413432
$ws.Range("K24").Value = 1.039814352689594
$ws.Range("L24").Value = 1.044110804079251
$ws.Range("M24").Value = 1.055618638047911
$ws.Range("N24").Value = 1.016324961249691
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032188262109183
$ws.Range("D25").Value = 1.037495660136872
$ws.Range("E25").Value = 1.042038569465596
$ws.Range("F25").Value = 1.0538704577423
$ws.Range("I25").Value = 1.039751662816409
$ws.Range("J25").Value = 1.037802899723618
$ws.Range("K25").Value = 1.040542869960728
$ws.Range("L25").Value = 1.045071620911137
$ws.Range("M25").Value = 1.056867259020863
$ws.Range("N25").Value = 1.016579533160116
